$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 165, shifting the existing
# rows 165-247 down to 167-249 (matches the dimension growing from
# A1:T247 to A1:T249).
$ws.Rows("165:166").Insert()

# New row 165: Crimpson Seedless / Especial
$ws.Cells.Item(165, 1).Value = 7
$ws.Cells.Item(165, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(165, 3).Value = "Ñuble"
$ws.Cells.Item(165, 4).Value = 45089
$ws.Cells.Item(165, 5).Value = 16
$ws.Cells.Item(165, 6).Value = "Fruta"
$ws.Cells.Item(165, 7).Value = 100109
$ws.Cells.Item(165, 8).Value = "Uva"
$ws.Cells.Item(165, 9).Value = 100109001
$ws.Cells.Item(165, 10).Value = "Uva"
$ws.Cells.Item(165, 11).Value = "Crimpson Seedless"
$ws.Cells.Item(165, 12).Value = "Especial"
$ws.Cells.Item(165, 13).Value = 80
$ws.Cells.Item(165, 14).Value = 13000
$ws.Cells.Item(165, 15).Value = 13000
$ws.Cells.Item(165, 16).Value = 13000
$ws.Cells.Item(165, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(165, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(165, 19).Value = 722
$ws.Cells.Item(165, 20).Value = 18

# New row 166: Crimpson Seedless / Primera
$ws.Cells.Item(166, 1).Value = 7
$ws.Cells.Item(166, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(166, 3).Value = "Ñuble"
$ws.Cells.Item(166, 4).Value = 45089
$ws.Cells.Item(166, 5).Value = 16
$ws.Cells.Item(166, 6).Value = "Fruta"
$ws.Cells.Item(166, 7).Value = 100109
$ws.Cells.Item(166, 8).Value = "Uva"
$ws.Cells.Item(166, 9).Value = 100109001
$ws.Cells.Item(166, 10).Value = "Uva"
$ws.Cells.Item(166, 11).Value = "Crimpson Seedless"
$ws.Cells.Item(166, 12).Value = "Primera"
$ws.Cells.Item(166, 13).Value = 50
$ws.Cells.Item(166, 14).Value = 12000
$ws.Cells.Item(166, 15).Value = 12000
$ws.Cells.Item(166, 16).Value = 12000
$ws.Cells.Item(166, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(166, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(166, 19).Value = 667
$ws.Cells.Item(166, 20).Value = 18
